# VSIG Trial Balance - second page
# Commit: "Changed all input from notes to smallcaps. -> update categories not
# done yet. still in progress"
#
# This adds a title block (company name / report title / report period) in
# column E of the header rows, and a "category" column (H) next to every
# account line that maps the raw TB account name to a standardised category
# name. It also totals the new category column's check (G56-F56) in H56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Title block (column E, rows 2-8) --------------------------------
$ws.Range("E2").Value = "VSIG Pte. Ltd."
$ws.Range("E3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("E7").Value = "Trial Balance"
$ws.Range("E8").Value = "December 2015"

# ---- Category column (column H, rows 11-54) ---------------------------
# Bank balances
$ws.Range("H11").Value = "Bank Balances"
$ws.Range("H12").Value = "Bank Balances"
$ws.Range("H13").Value = "Bank Balances"

# Trade receivables
$ws.Range("H14").Value = "Trade Receivables"
$ws.Range("H15").Value = "Trade Receivables"

# Plant and equipment
$ws.Range("H16").Value = "Plant and Equipment"
$ws.Range("H17").Value = "Plant and Equipment"
$ws.Range("H18").Value = "Plant and Equipment"
$ws.Range("H19").Value = "Plant and Equipment"

# Deposits
$ws.Range("H20").Value = "Deposits"

# Prepayments
$ws.Range("H21").Value = "Prepayments"

# Trade payables
$ws.Range("H22").Value = "Trade Payables"
$ws.Range("H23").Value = "Trade Payables"

# (H24, H25 - GST Collected / GST Paid - left uncategorised)

# GST payables
$ws.Range("H26").Value = "GST Payables"

# Accruals
$ws.Range("H27").Value = "Accruals"

# Amount owing to a shareholder
$ws.Range("H28").Value = "Amount owing to a Shareholder"

# Income tax payables
$ws.Range("H29").Value = "Income Tax Payables"

# Share capital
$ws.Range("H30").Value = "Share Capital"

# Retained profits
$ws.Range("H31").Value = "Retained Profits"

# Revenue
$ws.Range("H32").Value = "Revenue"

# Cost of sales
$ws.Range("H33").Value = "Cost of Sales"

# Operating expenses
$ws.Range("H34").Value = "Accounting Fee"
$ws.Range("H35").Value = "Administrative Expenses"
$ws.Range("H36").Value = "Bank Charges"
$ws.Range("H37").Value = "Compilation Fee"
$ws.Range("H38").Value = "Depreciation"
$ws.Range("H39").Value = "Entertainment"
$ws.Range("H40").Value = "Freight Charges"
$ws.Range("H41").Value = "Internet Expenses"
$ws.Range("H42").Value = "Late Penalty"
$ws.Range("H43").Value = "Nominee Director Fee"
$ws.Range("H44").Value = "Office Supplies"
$ws.Range("H45").Value = "Postage and Courier"
$ws.Range("H46").Value = "Professional Fee"
$ws.Range("H47").Value = "Secretarial Fee"
$ws.Range("H48").Value = "Taxation Fee"
$ws.Range("H49").Value = "Telephone Expenses"
$ws.Range("H50").Value = "Salaries"
$ws.Range("H51").Value = "Skill Development Levy & SINDA"

# Exchange gains
$ws.Range("H52").Value = "Exchange Gain - Trade"
$ws.Range("H53").Value = "Exchange Gain - Non-trade"

# Income tax expense
$ws.Range("H54").Value = "Income Tax Expense"

# ---- Check total for the new category column --------------------------
$ws.Range("H56").Formula = "=G56-F56"

# ---- Restore the cursor to where the author last left off -------------
$ws.Range("D24").Select()
